# Apply the diff:
#  - B256 changes from text "2001" to a real number 2001
#  - Three new rows (257, 258, 259) are appended with Title / Year / Link
#  - Row 259's Year (B259) must stay TEXT "2025" (not a number)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) B256: text "2001" -> numeric 2001
$ws.Range("B256").Value = 2001

# 2) Row 257: Nyan~ Neko Sugar Girls / 2010 / link
$ws.Range("A257").Value = "Nyan~ Neko Sugar Girls"
$ws.Range("B257").Value = 2010
$ws.Range("C257").Value = "https://letterboxd.com/film/nyan-neko-sugar-girls-2010/"

# 3) Row 258: Andor: A Disney+ Day Special Look / 2022 / link
$ws.Range("A258").Value = "Andor: A Disney+ Day Special Look"
$ws.Range("B258").Value = 2022
$ws.Range("C258").Value = "https://letterboxd.com/film/andor-a-disney-day-special-look/"

# 4) Row 259: Vanisher, Horizon Scraper / "2025" (kept as TEXT) / link
$ws.Range("A259").Value = "Vanisher, Horizon Scraper"
$ws.Range("B259").NumberFormat = "@"
$ws.Range("B259").Value = "2025"
$ws.Range("C259").Value = "https://letterboxd.com/film/vanisher-horizon-scraper/"
